# Generate Report for Archive
#
# The underlying e2e report data got re-sorted: the row that used to be the
# 3rd data row (6b7eca85-...) moved up to become the 1st data row, and the
# two rows that used to be 1st/2nd (775e6aa4-..., 7dc4498a-...) each shifted
# down by one row. The 4th data row (81015c87-...) stays put.
#
# This script re-creates that reordering on the "Overview", "zh-cn" and
# "de-de" worksheets, and fixes up the hyperlink display text to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A,B,E,F,G carry per-file data (rows 2-4 move)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "6b7eca85-f4a3-4bcf-a4ff-722720788659.md"
$ws.Range("B2").Value = "e2e\6b7eca85-f4a3-4bcf-a4ff-722720788659.md"
$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"
$ws.Range("G2").Value = "2016-10-26 07:13:31"

$ws.Range("A3").Value = "775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.md"
$ws.Range("B3").Value = "e2e\775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.md"
$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"
$ws.Range("G3").Value = "2016-10-26 07:12:09"

$ws.Range("A4").Value = "7dc4498a-563c-414d-9ca2-ef5828de0707.md"
$ws.Range("B4").Value = "e2e\7dc4498a-563c-414d-9ca2-ef5828de0707.md"
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"
$ws.Range("G4").Value = "2016-10-26 07:12:09"

# Hyperlinks: positions keep pointing at the same targets (rId2/3/4/5), but
# the visible display text must follow the file that now sits in that row.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97b21b0da07e90d7439ad8308bb0d882cd3f6104/e2e/775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.md", "", "", "e2e\6b7eca85-f4a3-4bcf-a4ff-722720788659.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97b21b0da07e90d7439ad8308bb0d882cd3f6104/e2e/7dc4498a-563c-414d-9ca2-ef5828de0707.md", "", "", "e2e\775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f7caa3e7ddc1c3ded62efe71b2102fe232b0dc5/e2e/6b7eca85-f4a3-4bcf-a4ff-722720788659.md", "", "", "e2e\7dc4498a-563c-414d-9ca2-ef5828de0707.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25822f81c065563471e462149b8adef7da343fb9/e2e/81015c87-6f0f-49eb-bec6-73ae877e7b2e.md", "", "", "e2e\81015c87-6f0f-49eb-bec6-73ae877e7b2e.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A,C,G,H carry per-file data (rows 2-4 move)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "6b7eca85-f4a3-4bcf-a4ff-722720788659.md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("G2").Value = "6b7eca85-f4a3-4bcf-a4ff-722720788659.c0f883fac2efdaf8506e21f31255ff67ec8794bd.zh-cn.xlf"
$ws.Range("H2").Value = "2016-10-26 07:13:20"

$ws.Range("A3").Value = "775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("G3").Value = "775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.df7c5697ef593da6b93e1585e394ece6c80fb5d4.zh-cn.xlf"
$ws.Range("H3").Value = "2016-10-26 07:11:52"

$ws.Range("A4").Value = "7dc4498a-563c-414d-9ca2-ef5828de0707.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "7dc4498a-563c-414d-9ca2-ef5828de0707.2e05772bd206aeb3a6f5e4927e454d2bdf2fe46e.zh-cn.xlf"
$ws.Range("H4").Value = "2016-10-26 07:11:52"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97b21b0da07e90d7439ad8308bb0d882cd3f6104/e2e/775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.md", "", "", "6b7eca85-f4a3-4bcf-a4ff-722720788659.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97b21b0da07e90d7439ad8308bb0d882cd3f6104/e2e/7dc4498a-563c-414d-9ca2-ef5828de0707.md", "", "", "775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f7caa3e7ddc1c3ded62efe71b2102fe232b0dc5/e2e/6b7eca85-f4a3-4bcf-a4ff-722720788659.md", "", "", "7dc4498a-563c-414d-9ca2-ef5828de0707.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25822f81c065563471e462149b8adef7da343fb9/e2e/81015c87-6f0f-49eb-bec6-73ae877e7b2e.md", "", "", "81015c87-6f0f-49eb-bec6-73ae877e7b2e.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": columns A,C,G,H carry per-file data (rows 2-4 move)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "6b7eca85-f4a3-4bcf-a4ff-722720788659.md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("G2").Value = "6b7eca85-f4a3-4bcf-a4ff-722720788659.c0f883fac2efdaf8506e21f31255ff67ec8794bd.de-de.xlf"
$ws.Range("H2").Value = "2016-10-26 07:13:31"

$ws.Range("A3").Value = "775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("G3").Value = "775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.df7c5697ef593da6b93e1585e394ece6c80fb5d4.de-de.xlf"
$ws.Range("H3").Value = "2016-10-26 07:12:09"

$ws.Range("A4").Value = "7dc4498a-563c-414d-9ca2-ef5828de0707.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "7dc4498a-563c-414d-9ca2-ef5828de0707.2e05772bd206aeb3a6f5e4927e454d2bdf2fe46e.de-de.xlf"
$ws.Range("H4").Value = "2016-10-26 07:12:09"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97b21b0da07e90d7439ad8308bb0d882cd3f6104/e2e/775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.md", "", "", "6b7eca85-f4a3-4bcf-a4ff-722720788659.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97b21b0da07e90d7439ad8308bb0d882cd3f6104/e2e/7dc4498a-563c-414d-9ca2-ef5828de0707.md", "", "", "775e6aa4-f2c8-4e0d-80ef-2ebbf21e93a2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f7caa3e7ddc1c3ded62efe71b2102fe232b0dc5/e2e/6b7eca85-f4a3-4bcf-a4ff-722720788659.md", "", "", "7dc4498a-563c-414d-9ca2-ef5828de0707.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25822f81c065563471e462149b8adef7da343fb9/e2e/81015c87-6f0f-49eb-bec6-73ae877e7b2e.md", "", "", "81015c87-6f0f-49eb-bec6-73ae877e7b2e.md") | Out-Null
